{"js": "// Update the worksheet date header and the 25 division-problem answers\n// in the single table, in document order. Source data does contain\n// duplicate \"problem\" strings (e.g. \"75\u00f79=8, 3\" appears twice with two\n// different replacements), so replacement MUST be positional (row/col\n// indexed), not a global find-and-replace.\n\nconst body = context.document.body;\n\n// 1) Date paragraph above the table: \"2026-02-06 Friday\" -> \"2026-02-07 Saturday\"\nconst firstParagraph = body.paragraphs.getFirst();\nfirstParagraph.load(\"text\");\nawait context.sync();\n\nif (firstParagraph.text.trim() === \"2026-02-06 Friday\") {\n  firstParagraph.getRange().insertText(\"2026-02-07 Saturday\", \"Replace\");\n}\n\n// 2) Table of division problems: 5 populated rows x 5 columns (rows 0,4,8,12,16\n// of the 20-row table; the other rows are blank spacer rows).\nconst table = body.tables.getFirst();\n\nconst newValues = [\n  [\"91\u00f76=15, 1\", \"82\u00f75=16, 2\", \"55\u00f73=18, 1\", \"44\u00f76=7, 2\", \"80\u00f72=40, 0\"],\n  [\"16\u00f75=3, 1\", \"25\u00f73=8, 1\", \"41\u00f76=6, 5\", \"20\u00f78=2, 4\", \"39\u00f73=13, 0\"],\n  [\"54\u00f72=27, 0\", \"34\u00f77=4, 6\", \"65\u00f76=10, 5\", \"68\u00f74=17, 0\", \"31\u00f74=7, 3\"],\n  [\"79\u00f75=15, 4\", \"20\u00f79=2, 2\", \"51\u00f73=17, 0\", \"91\u00f77=13, 0\", \"51\u00f77=7, 2\"],\n  [\"26\u00f77=3, 5\", \"50\u00f77=7, 1\", \"90\u00f73=30, 0\", \"44\u00f74=11, 0\", \"36\u00f74=9, 0\"],\n];\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const row = dataRowIndexes[i];\n  for (let col = 0; col < newValues[i].length; col++) {\n    table.getCell(row, col).value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 division-problem answers\n# in the single table, in document order. Source data contains duplicate\n# \"problem\" strings (e.g. \"75\u00f79=8, 3\" appears twice with two different\n# replacements), so replacement MUST be positional (row/col indexed via\n# Table.Cell(r, c)), not a global Find/Replace.\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph above the table.\n$p = $d.Paragraphs(1)\nif ($p.Range.Text.TrimEnd(\"`r\") -eq \"2026-02-06 Friday\") {\n    $p.Range.Text = \"2026-02-07 Saturday\"\n}\n\n# 2) Table of division problems: 5 populated rows x 5 columns (table rows\n# 1, 5, 9, 13, 17 in 1-based COM numbering; the rows in between are blank\n# spacer rows).\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    ,@(\"91\u00f76=15, 1\", \"82\u00f75=16, 2\", \"55\u00f73=18, 1\", \"44\u00f76=7, 2\", \"80\u00f72=40, 0\")\n    ,@(\"16\u00f75=3, 1\", \"25\u00f73=8, 1\", \"41\u00f76=6, 5\", \"20\u00f78=2, 4\", \"39\u00f73=13, 0\")\n    ,@(\"54\u00f72=27, 0\", \"34\u00f77=4, 6\", \"65\u00f76=10, 5\", \"68\u00f74=17, 0\", \"31\u00f74=7, 3\")\n    ,@(\"79\u00f75=15, 4\", \"20\u00f79=2, 2\", \"51\u00f73=17, 0\", \"91\u00f77=13, 0\", \"51\u00f77=7, 2\")\n    ,@(\"26\u00f77=3, 5\", \"50\u00f77=7, 1\", \"90\u00f73=30, 0\", \"44\u00f74=11, 0\", \"36\u00f74=9, 0\")\n)\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $row = $dataRows[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
